# Insert a new weekly price record at row 146 for
# "Terminal Hortofrutícola Agro Chillán - Pepino ensalada", pushing the
# existing rows 146:171 down to 147:172.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 146 and below down by one row.
$ws.Rows("146:146").Insert()

# Populate the newly inserted row with the new data point.
$ws.Range("A146").Value = 7
$ws.Range("B146").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C146").Value = "Ñuble"
$ws.Range("D146").Value = 44511
$ws.Range("E146").Value = 16
$ws.Range("F146").Value = 100112043
$ws.Range("G146").Value = "Pepino ensalada"
$ws.Range("H146").Value = "Sin especificar"
$ws.Range("I146").Value = "Primera"
$ws.Range("J146").Value = 100
$ws.Range("K146").Value = 8000
$ws.Range("L146").Value = 9000
$ws.Range("M146").Value = 8500
$ws.Range("N146").Value = "$/caja 80 unidades"
$ws.Range("O146").Value = "Región del Maule"
$ws.Range("P146").Value = 106
$ws.Range("Q146").Value = 80
$ws.Range("R146").Value = "Hortaliza"
